$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.077.07"
$ws.Range("E2").Value = "  +3.70%  "

# Row 3
$ws.Range("D3").Value = "1.727.67"
$ws.Range("E3").Value = "  +2.94%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.09"
$ws.Range("E5").Value = "  +1.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.524"
$ws.Range("E6").Value = "  +1.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.13"
$ws.Range("E8").Value = "  +13.52%  "

# Row 9
$ws.Range("E9").Value = "  +3.37%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0633"
$ws.Range("E10").Value = "  +1.84%  "

# Row 11
$ws.Range("E11").Value = "  +1.61%  "

# Row 12
$ws.Range("D12").Value = "1.970.38"
$ws.Range("E12").Value = "  +2.90%  "

# Row 13
$ws.Range("D13").Value = "1.730.06"
$ws.Range("E13").Value = "  +2.89%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.28"
$ws.Range("E14").Value = "  +3.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.53"
$ws.Range("E16").Value = "  +2.22%  "

# Row 17
$ws.Range("D17").Value = "28.039.49"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.29"
$ws.Range("E18").Value = "  +2.68%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0756"
$ws.Range("E19").Value = "  +1.91%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.90"
$ws.Range("E20").Value = "  -2.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.64"
$ws.Range("E22").Value = "  +3.73%  "

# Row 23
$ws.Range("E23").Value = "  +4.11%  "

# Row 24
$ws.Range("E24").Value = "  +0.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.95"
$ws.Range("E25").Value = "  +1.64%  "

# Row 26
$ws.Range("E26").Value = "  +4.17%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.72"

# Row 28
$ws.Range("E28").Value = "  +1.49%  "

# Row 29
$ws.Range("E29").Value = "  -0.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0511"
$ws.Range("E30").Value = "  +2.50%  "

# Row 31
$ws.Range("E31").Value = "  +2.26%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.45"

# Row 33
$ws.Range("D33").Value = "1.493.95"
$ws.Range("E33").Value = "  -3.89%  "

# Row 34
$ws.Range("E34").Value = "  +2.68%  "

# Row 35
$ws.Range("E35").Value = "  -1.75%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.609"
$ws.Range("E36").Value = "  +1.40%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.954"
$ws.Range("E37").Value = "  +3.03%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.41"
$ws.Range("E38").Value = "  +1.11%  "

# Row 39
$ws.Range("E39").Value = "  +0.44%  "

# Row 40
$ws.Range("E40").Value = "  +0.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "70.72"
$ws.Range("E41").Value = "  +3.82%  "

# Row 42
$ws.Range("E42").Value = "  +3.92%  "

# Row 43
$ws.Range("E43").Value = "  -0.04%  "

# Row 44
$ws.Range("E44").Value = "  +2.31%  "

# Row 45
$ws.Range("D45").Value = "1.875.00"
$ws.Range("E45").Value = "  +2.75%  "

# Row 46
$ws.Range("E46").Value = "  +2.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.80"
$ws.Range("E47").Value = "  +14.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "91.17"
$ws.Range("E48").Value = "  +0.52%  "

# Row 49
$ws.Range("E49").Value = "  +4.44%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.24"
$ws.Range("E50").Value = "  +2.38%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").Value = "  +0.79%  "
